$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) - update "想去人数" (F) column values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 284
$ws1.Range("F3").Value = 1167
$ws1.Range("F4").Value = 16626
$ws1.Range("F5").Value = 20
$ws1.Range("F6").Value = 1631
$ws1.Range("F10").Value = 209
$ws1.Range("F12").Value = 11568
$ws1.Range("F13").Value = 25
$ws1.Range("F14").Value = 1246
$ws1.Range("F15").Value = 4575
$ws1.Range("F16").Value = 410
$ws1.Range("F17").Value = 401
$ws1.Range("F19").Value = 875
$ws1.Range("F20").Value = 332
$ws1.Range("F21").Value = 150

# Sheet "全部类型" (sheet4) - update "想去人数" (F) column values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 284
$ws4.Range("F4").Value = 1167
$ws4.Range("F5").Value = 16626
$ws4.Range("F6").Value = 20
$ws4.Range("F7").Value = 1631
$ws4.Range("F11").Value = 209
$ws4.Range("F15").Value = 11568
$ws4.Range("F16").Value = 25
$ws4.Range("F17").Value = 1247
$ws4.Range("F18").Value = 4575
$ws4.Range("F19").Value = 410
$ws4.Range("F20").Value = 401
$ws4.Range("F22").Value = 875
$ws4.Range("F23").Value = 332
$ws4.Range("F24").Value = 150
